# Daily attendance processing - 2026-01-30 00:06:55
# Rotates the "Recorded By" (column G) comma-separated list of names/emails
# left by one position (first entry moves to the end) for the specific rows
# that were touched during this processing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 4, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 28, 30, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 54, 56, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $current = [string]$cell.Value2
    $parts = $current -split ',\s*'
    if ($parts.Count -gt 1) {
        $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
        $cell.Value = [string]::Join(', ', $rotated)
    }
}
